$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.851.10"
$ws.Range("E2").Value = "  -2.68%  "
$ws.Range("D3").Value = "2.575.60"
$ws.Range("E3").Value = "  -5.19%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'552.74"
$ws.Range("E5").Value = "  -1.18%  "
$ws.Range("D6").Value = "'154.52"
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.599"
$ws.Range("E8").Value = "  +1.25%  "
$ws.Range("E9").Value = "  -2.63%  "
$ws.Range("D10").Value = "'0.163"
$ws.Range("E10").Value = "  -1.22%  "
$ws.Range("D12").Value = "'0.365"
$ws.Range("E12").Value = "  -1.93%  "
$ws.Range("D13").Value = "3.033.42"
$ws.Range("E13").Value = "  -5.13%  "
$ws.Range("D14").Value = "'25.53"
$ws.Range("E14").Value = "  -3.57%  "
$ws.Range("D15").Value = "61.764.10"
$ws.Range("E15").Value = "  -2.56%  "
$ws.Range("D16").Value = "'0.0000144"
$ws.Range("E16").Value = "  -2.28%  "
$ws.Range("D17").Value = "2.578.73"
$ws.Range("E17").Value = "  -5.16%  "
$ws.Range("D18").Value = "'11.59"
$ws.Range("E18").Value = "  -4.53%  "
$ws.Range("E19").Value = "  -2.13%  "
$ws.Range("D20").Value = "'338.13"
$ws.Range("E20").Value = "  -3.13%  "
$ws.Range("D21").Value = "'6.07"
$ws.Range("E21").Value = "  -5.22%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("E23").Value = "  -3.53%  "
$ws.Range("D24").Value = "'63.33"
$ws.Range("E24").Value = "  -1.22%  "
$ws.Range("E25").Value = "  -1.03%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("D27").Value = "'8.10"
$ws.Range("E27").Value = "  -0.76%  "
$ws.Range("D28").Value = "'7.41"
$ws.Range("E28").Value = "  +3.67%  "
$ws.Range("E29").Value = "  -4.45%  "
$ws.Range("E30").Value = "  -1.41%  "
$ws.Range("D31").Value = "'1.31"
$ws.Range("E31").Value = "  -3.87%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("D33").Value = "'159.35"
$ws.Range("E33").Value = "  -3.06%  "
$ws.Range("D34").Value = "'4.74"
$ws.Range("E34").Value = "  -1.66%  "
$ws.Range("D35").Value = "'19.22"
$ws.Range("E35").Value = "  -3.05%  "
$ws.Range("E36").Value = "  -3.45%  "
$ws.Range("E37").Value = "  +1.52%  "
$ws.Range("D38").Value = "'338.12"
$ws.Range("E38").Value = "  -2.99%  "
$ws.Range("D39").Value = "'0.939"
$ws.Range("E39").Value = "  -2.09%  "
$ws.Range("D40").Value = "'6.01"
$ws.Range("E40").Value = "  -1.13%  "
$ws.Range("D41").Value = "'3.97"
$ws.Range("E41").Value = "  -0.52%  "
$ws.Range("D42").Value = "'37.66"
$ws.Range("E42").Value = "  -1.43%  "
$ws.Range("D43").Value = "'20.65"
$ws.Range("E43").Value = "  -2.99%  "
$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").Value = "2.142.45"
$ws.Range("E45").Value = "  +1.67%  "
$ws.Range("D46").Value = "'0.606"
$ws.Range("E46").Value = "  -2.93%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "'0.0549"
$ws.Range("E47").Value = "  -3.98%  "
$ws.Range("D48").Value = "'19.67"
$ws.Range("E48").Value = "  -4.62%  "
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").Value = "'10.93"
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("D50").Value = "'0.0968"
$ws.Range("E50").Value = "  -1.41%  "
$ws.Range("E51").Value = "  -2.06%  "
